$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "sunny weather"
$ws.Range("B2").Value = "晴れ|はれ"
$ws.Range("A3").Value = "rain"
$ws.Range("B3").Value = "雨|あめ"
$ws.Range("A4").Value = "cloudy weather"
$ws.Range("B4").Value = "曇り|くもり"
$ws.Range("A5").Value = "snow"
$ws.Range("B5").Value = "雪|ゆき"
$ws.Range("A6").Value = "weather forecast"
$ws.Range("B6").Value = "天気予報|てんきよほう"
$ws.Range("A7").Value = "temperature (weather)"
$ws.Range("B7").Value = "気温|きおん"
$ws.Range("A8").Value = "summer"
$ws.Range("B8").Value = "夏|なつ"
$ws.Range("A9").Value = "winter"
$ws.Range("B9").Value = "冬|ふゆ"
$ws.Range("A10").Value = "this morning"
$ws.Range("B10").Value = "今朝|けさ"
$ws.Range("A11").Value = "the day after tomorrow"
$ws.Range("B11").Value = "あさって"
$ws.Range("A12").Value = "every week"
$ws.Range("B12").Value = "毎週|まいしゅう"
$ws.Range("A13").Value = "this month"
$ws.Range("B13").Value = "今月|こんげつ"
$ws.Range("A14").Value = "next month"
$ws.Range("B14").Value = "来月|らいげつ"
$ws.Range("A15").Value = "office worker"
$ws.Range("B15").Value = "会社員|かいしゃいん"
$ws.Range("A16").Value = "job; work; occupation"
$ws.Range("B16").Value = "仕事|しごと"
$ws.Range("A17").Value = "camera"
$ws.Range("B17").Value = "カメラ"
$ws.Range("A18").Value = "karaoke"
$ws.Range("B18").Value = "カラオケ"
$ws.Range("A19").Value = "place"
$ws.Range("B19").Value = "所|ところ"
$ws.Range("A20").Value = "tomato"
$ws.Range("B20").Value = "トマト"
$ws.Range("A21").Value = "chopsticks"
$ws.Range("B21").Value = "はし"
$ws.Range("A22").Value = "party"
$ws.Range("B22").Value = "パーティー"
$ws.Range("A23").Value = "barbecue"
$ws.Range("B23").Value = "バーベキュー"
$ws.Range("A24").Value = "homestay; living with a local family"
$ws.Range("B24").Value = "ホームステイ"
$ws.Range("A25").Value = "bath"
$ws.Range("B25").Value = "お風呂|おふろ"
$ws.Range("A26").Value = "Spain"
$ws.Range("B26").Value = "スペイン"
$ws.Range("A27").Value = "something"
$ws.Range("B27").Value = "何か|なにか"
$ws.Range("A28").Value = "skillful; good at (～が)"
$ws.Range("B28").Value = "上手|じょうず(な)"
$ws.Range("A29").Value = "clumsy; poor at (～が)"
$ws.Range("B29").Value = "下手|へた(な)"
$ws.Range("A30").Value = "famous"
$ws.Range("B30").Value = "有名|ゆうめい(な)"
$ws.Range("A31").Value = "to wash (～を)"
$ws.Range("B31").Value = "洗う|あらう"
$ws.Range("A32").Value = "to say"
$ws.Range("B32").Value = "言う|いう"
$ws.Range("A33").Value = "to need (～が)"
$ws.Range("B33").Value = "いる"
$ws.Range("A34").Value = "to be late"
$ws.Range("B34").Value = "遅くなる|おそくなる"
$ws.Range("A35").Value = "to take a bath"
$ws.Range("B35").Value = "お風呂に入る|おふろにはいる"
$ws.Range("A36").Value = "to think"
$ws.Range("B36").Value = "思う|おもう"
$ws.Range("A37").Value = "to cut (～を)"
$ws.Range("B37").Value = "切る|きる"
$ws.Range("A38").Value = "to make (～を)"
$ws.Range("B38").Value = "作る|つくる"
$ws.Range("A39").Value = "rain falls"
$ws.Range("B39").Value = "雨が降る|あめがふる"
$ws.Range("A40").Value = "snow falls"
$ws.Range("B40").Value = "雪が降る|ゆきがふる"
$ws.Range("A41").Value = "to take (a thing) (～を)"
$ws.Range("B41").Value = "持っていく|もっていく"
$ws.Range("A42").Value = "to throw away (～を)"
$ws.Range("B42").Value = "捨てる|すてる"
$ws.Range("A43").Value = "to begin (～を)"
$ws.Range("B43").Value = "始める|はじめる"
$ws.Range("A44").Value = "to drive (～を)"
$ws.Range("B44").Value = "運転する|うんてんする"
$ws.Range("A45").Value = "to do laundry (～を)"
$ws.Range("B45").Value = "洗濯する|せんたくする"
$ws.Range("A46").Value = "to clean (～を)"
$ws.Range("B46").Value = "掃除する|そうじする"
$ws.Range("A47").Value = "to cook"
$ws.Range("B47").Value = "料理する|りょうりする"
$ws.Range("A48").Value = "uh-huh; yes"
$ws.Range("B48").Value = "うん"
$ws.Range("A49").Value = "uh-uh; no"
$ws.Range("B49").Value = "ううん"
$ws.Range("A50").Value = "always"
$ws.Range("B50").Value = "いつも"
$ws.Range("A51").Value = "(do something) late"
$ws.Range("B51").Value = "遅く|おそく"
$ws.Range("A52").Value = "Cheers! (a toast)"
$ws.Range("B52").Value = "乾杯|かんぱい"
$ws.Range("A53").Value = "all (of the people) together"
$ws.Range("B53").Value = "みんなで"
$ws.Range("A54").Value = "That's too bad."
$ws.Range("B54").Value = "残念(ですね)|ざんねん(ですね)"
$ws.Range("A55").Value = "not...yet (w/negative)"
$ws.Range("B55").Value = "まだ"
$ws.Range("A56").Value = "about...; concerning..."
$ws.Range("B56").Value = "～について"
$ws.Range("A57").Value = "...degrees (temperature)"
$ws.Range("B57").Value = "～度|～ど"
$ws.Range("A58").Value = "how"
$ws.Range("B58").Value = "どう"
$ws.Range("A59").Value = "rice"
$ws.Range("B59").Value = "ご飯|ごはん"
$ws.Range("A60").Value = "side dish"
$ws.Range("B60").Value = "おかず"
$ws.Range("A61").Value = "miso soup"
$ws.Range("B61").Value = "みそ汁|みそしる"
$ws.Range("A62").Value = "set meal"
$ws.Range("B62").Value = "定食|ていしょく"
$ws.Range("A63").Value = "curry with rice"
$ws.Range("B63").Value = "カレーライス"
$ws.Range("A64").Value = "rice balls"
$ws.Range("B64").Value = "おにぎり"
$ws.Range("A65").Value = "ramen noodles"
$ws.Range("B65").Value = "ラーメン"
$ws.Range("A66").Value = "udon noodles"
$ws.Range("B66").Value = "うどん"
$ws.Range("A67").Value = "pasta"
$ws.Range("B67").Value = "パスタ"
$ws.Range("A68").Value = "dumplings"
$ws.Range("B68").Value = "ぎょうざ"
$ws.Range("A69").Value = "beef rice bowl"
$ws.Range("B69").Value = "牛丼|ぎゅうどん"
$ws.Range("A70").Value = "hamburger steak"
$ws.Range("B70").Value = "ハンバーグ"
$ws.Range("A71").Value = "raw seafood"
$ws.Range("B71").Value = "さしみ"
$ws.Range("A72").Value = "savory pancake"
$ws.Range("B72").Value = "お好み焼き|おこのみやき"
$ws.Range("A73").Value = "toast"
$ws.Range("B73").Value = "トースト"
$ws.Range("A74").Value = "soup"
$ws.Range("B74").Value = "スープ"
$ws.Range("A75").Value = "yogurt"
$ws.Range("B75").Value = "ヨーグルト"
$ws.Range("A76").Value = "broiled fish"
$ws.Range("B76").Value = "焼き魚|やきざかな"
$ws.Range("A77").Value = "egg"
$ws.Range("B77").Value = "たまご"
$ws.Range("A78").Value = "Do you have avocado?"
$ws.Range("B78").Value = "アボカドはありますか。"
$ws.Range("A79").Value = "Is there alcohol in this?"
$ws.Range("B79").Value = "この中にお酒が入っていますか。|このなかにおさけがはいっていますか。"
$ws.Range("A80").Value = "Is this halal?"
$ws.Range("B80").Value = "これはハラルフードですか。"
$ws.Range("A81").Value = "I have an allergy to peanuts."
$ws.Range("B81").Value = "ピーナッツアレルギーがあります。"
$ws.Range("A82").Value = "Can I have a bag?"
$ws.Range("B82").Value = "袋をお願いします。|ふくろをおねがいします。"
$ws.Range("A83").Value = "fruit"
$ws.Range("B83").Value = "果物|くだもの"
$ws.Range("A84").Value = "strawberry"
$ws.Range("B84").Value = "いちご"
$ws.Range("A85").Value = "watermelon"
$ws.Range("B85").Value = "すいか"
$ws.Range("A86").Value = "mandarin orange"
$ws.Range("B86").Value = "みかん"
$ws.Range("A87").Value = "apple"
$ws.Range("B87").Value = "りんご"
$ws.Range("A88").Value = "peach"
$ws.Range("B88").Value = "もも"
$ws.Range("A89").Value = "grapes"
$ws.Range("B89").Value = "ぶどう"
$ws.Range("A90").Value = "vegetable"
$ws.Range("B90").Value = "野菜|やさい"
$ws.Range("A91").Value = "carrot"
$ws.Range("B91").Value = "にんじん"
$ws.Range("A92").Value = "onion"
$ws.Range("B92").Value = "たまねぎ"
$ws.Range("A93").Value = "potato"
$ws.Range("B93").Value = "じゃがいも"
$ws.Range("A94").Value = "eggplant"
$ws.Range("B94").Value = "なす"
$ws.Range("A95").Value = "cucumber"
$ws.Range("B95").Value = "きゅうり"
$ws.Range("A96").Value = "cabbage"
$ws.Range("B96").Value = "キャベツ"
$ws.Range("A97").Value = "meat"
$ws.Range("B97").Value = "肉|にく"
$ws.Range("A98").Value = "beef"
$ws.Range("B98").Value = "牛肉|ぎゅうにく"
$ws.Range("A99").Value = "pork"
$ws.Range("B99").Value = "豚肉|ぶたにく"
$ws.Range("A100").Value = "chicken"
$ws.Range("B100").Value = "鶏肉|とりにく"
$ws.Range("A101").Value = "office worker"
$ws.Range("B101").Value = "会社員|かいしゃいん"
$ws.Range("A102").Value = "store clerk"
$ws.Range("B102").Value = "店員|てんいん"
$ws.Range("A103").Value = "member"
$ws.Range("B103").Value = "会員|かいいん"
$ws.Range("A104").Value = "station attendant"
$ws.Range("B104").Value = "駅員|えきいん"
$ws.Range("A105").Value = "new"
$ws.Range("B105").Value = "新しい|あたらしい"
$ws.Range("A106").Value = "newspaper"
$ws.Range("B106").Value = "新聞|しんぶん"
$ws.Range("A107").Value = "Bullet Train"
$ws.Range("B107").Value = "新幹線|しんかんせん"
$ws.Range("A108").Value = "fresh"
$ws.Range("B108").Value = "新鮮な|しんせんな"
$ws.Range("A109").Value = "to listen"
$ws.Range("B109").Value = "聞く|きく"
$ws.Range("A110").Value = "can be heard"
$ws.Range("B110").Value = "聞こえる|きこえる"
$ws.Range("A111").Value = "to make"
$ws.Range("B111").Value = "作る|つくる"
$ws.Range("A112").Value = "composition"
$ws.Range("B112").Value = "作文|さくぶん"
$ws.Range("A113").Value = "artistic piece"
$ws.Range("B113").Value = "作品|さくひん"
$ws.Range("A114").Value = "author"
$ws.Range("B114").Value = "作者|さくしゃ"
$ws.Range("A115").Value = "job"
$ws.Range("B115").Value = "仕事|しごと"
$ws.Range("A116").Value = "revenge"
$ws.Range("B116").Value = "仕返し|しかえし"
$ws.Range("A117").Value = "to serve; to work under"
$ws.Range("B117").Value = "仕える|つかえる"
$ws.Range("A118").Value = "thing"
$ws.Range("B118").Value = "事|こと"
$ws.Range("A119").Value = "fire"
$ws.Range("B119").Value = "火事|かじ"
$ws.Range("A120").Value = "meal"
$ws.Range("B120").Value = "食事|しょくじ"
$ws.Range("A121").Value = "reply"
$ws.Range("B121").Value = "返事|へんじ"
$ws.Range("A122").Value = "train"
$ws.Range("B122").Value = "電車|でんしゃ"
$ws.Range("A123").Value = "electricity; light"
$ws.Range("B123").Value = "電気|でんき"
$ws.Range("A124").Value = "telephone"
$ws.Range("B124").Value = "電話|でんわ"
$ws.Range("A125").Value = "battery"
$ws.Range("B125").Value = "電池|でんち"
$ws.Range("A126").Value = "microwave oven"
$ws.Range("B126").Value = "電子レンジ|でんしレンジ"
$ws.Range("A127").Value = "car"
$ws.Range("B127").Value = "車|くるま"
$ws.Range("A128").Value = "bicycle"
$ws.Range("B128").Value = "自転車|じてんしゃ"
$ws.Range("A129").Value = "wheelchair"
$ws.Range("B129").Value = "車いす|くるまいす"
$ws.Range("A130").Value = "parking lot"
$ws.Range("B130").Value = "駐車場|ちゅうしゃじょう"
$ws.Range("A131").Value = "to be absent; to rest"
$ws.Range("B131").Value = "休む|やすむ"
$ws.Range("A132").Value = "holiday; absence"
$ws.Range("B132").Value = "休み|やすみ"
$ws.Range("A133").Value = "holiday"
$ws.Range("B133").Value = "休日|きゅうじつ"
$ws.Range("A134").Value = "to say"
$ws.Range("B134").Value = "言う|いう"
$ws.Range("A135").Value = "linguistics"
$ws.Range("B135").Value = "言語学|げんごがく"
$ws.Range("A136").Value = "dialect"
$ws.Range("B136").Value = "方言|ほうげん"
$ws.Range("A137").Value = "word; language"
$ws.Range("B137").Value = "言葉|ことば"
$ws.Range("A138").Value = "to read"
$ws.Range("B138").Value = "読む|よむ"
$ws.Range("A139").Value = "reading books"
$ws.Range("B139").Value = "読書|どくしょ"
$ws.Range("A140").Value = "reading matter"
$ws.Range("B140").Value = "読み物|よみもの"
$ws.Range("A141").Value = "to think"
$ws.Range("B141").Value = "思う|おもう"
$ws.Range("A142").Value = "mysterious"
$ws.Range("B142").Value = "不思議な|ふしぎな"
$ws.Range("A143").Value = "to recall; to remember"
$ws.Range("B143").Value = "思い出す|おもいだす"
$ws.Range("A144").Value = "next"
$ws.Range("B144").Value = "次|つぎ"
$ws.Range("A145").Value = "second daughter"
$ws.Range("B145").Value = "次女|じじょ"
$ws.Range("A146").Value = "table of contents"
$ws.Range("B146").Value = "目次|もくじ"
$ws.Range("A147").Value = "next time"
$ws.Range("B147").Value = "次回|じかい"
$ws.Range("A148").Value = "what"
$ws.Range("B148").Value = "何|なに"
$ws.Range("A149").Value = "what time"
$ws.Range("B149").Value = "何時|なんじ"
$ws.Range("A150").Value = "how many people"
$ws.Range("B150").Value = "何人|なんにん"
$ws.Range("A151").Value = "something"
$ws.Range("B151").Value = "何か|なにか"
